$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $range1 = $ws.Range("B$row1" + ":AC$row1")
    $range2 = $ws.Range("B$row2" + ":AC$row2")
    $vals1 = $range1.Value()
    $vals2 = $range2.Value()
    $range1.Value = $vals2
    $range2.Value = $vals1
}

Swap-Rows 39 40
Swap-Rows 111 112
